# "Generate Report for Handback"
# The localization status report gets updated once a handback (zh-cn + de-de)
# comes back in sync with en-US: the per-language sheets grow two new
# columns (Latest Target File / Latest Handback File), the Status text
# flips from "Ready for handoff" to "Handed back: in sync with en-US"
# (which also flows through to the Overview roll-up), and the Latest
# Handback DateTime cells get populated with real timestamps.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdFile1  = "8313d808-0eab-4332-ba12-098343be3408.md"
$mdFile2  = "ef5b968b-20a2-4d60-8071-f3d595891d91.md"

$xlfZh1 = "8313d808-0eab-4332-ba12-098343be3408.dfc5d76faf8e29f36adae65bb55a559e724dda54.zh-cn.xlf"
$xlfZh2 = "ef5b968b-20a2-4d60-8071-f3d595891d91.66244de5adb47e7b3f8ee07635627c129c8fc3e4.zh-cn.xlf"
$xlfDe1 = "8313d808-0eab-4332-ba12-098343be3408.dfc5d76faf8e29f36adae65bb55a559e724dda54.de-de.xlf"
$xlfDe2 = "ef5b968b-20a2-4d60-8071-f3d595891d91.66244de5adb47e7b3f8ee07635627c129c8fc3e4.de-de.xlf"

# -- Overview sheet: the per-language status columns (B = zh-cn, C = de-de)
#    echo the same status text, so they flip too now that both languages
#    have been handed back.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# -- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column (C) flips for both rows.
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# New "Latest Target File" (F) / "Latest Handback File" (G) columns,
# mirroring the source markdown file / handoff xlf file respectively,
# each carrying its own hyperlink just like columns A and D do.
$wsZh.Range("F2").Value = $mdFile1
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e615ae7016605deeac12253686483921c023ef4/e2e/$mdFile1", "", "", $mdFile1) | Out-Null

$wsZh.Range("G2").Value = $xlfZh1
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/67537d25bf2ac946fba2a4af337f04ebb2e0031e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/$xlfZh1", "", "", $xlfZh1) | Out-Null

$wsZh.Range("F3").Value = $mdFile2
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e615ae7016605deeac12253686483921c023ef4/e2e/$mdFile2", "", "", $mdFile2) | Out-Null

$wsZh.Range("G3").Value = $xlfZh2
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/67537d25bf2ac946fba2a4af337f04ebb2e0031e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/$xlfZh2", "", "", $xlfZh2) | Out-Null

# Latest Handback DateTime (H) now has a real timestamp instead of the
# "never handed back" placeholder.
$wsZh.Range("H2").Value = "2016-03-25 07:56:24"
$wsZh.Range("H3").Value = "2016-03-25 07:56:24"

# -- de-de sheet ----------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("F2").Value = $mdFile1
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/2e615ae7016605deeac12253686483921c023ef4/e2e/$mdFile1", "", "", $mdFile1) | Out-Null

$wsDe.Range("G2").Value = $xlfDe1
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fe42363582669f540887a9fe8c4b18cfc5e283eb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/$xlfDe1", "", "", $xlfDe1) | Out-Null

$wsDe.Range("F3").Value = $mdFile2
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e615ae7016605deeac12253686483921c023ef4/e2e/$mdFile2", "", "", $mdFile2) | Out-Null

$wsDe.Range("G3").Value = $xlfDe2
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fe42363582669f540887a9fe8c4b18cfc5e283eb/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/$xlfDe2", "", "", $xlfDe2) | Out-Null

# de-de's handback happened a few seconds after zh-cn's.
$wsDe.Range("H2").Value = "2016-03-25 07:56:31"
$wsDe.Range("H3").Value = "2016-03-25 07:56:31"
